# Actualizacion 16 de Abril de 2020.
# Adds the new daily row (row 45, fecha 15-04-2020) to both sheets and
# updates the window scroll position / selection to match.

$wb = $excel.ActiveWorkbook

# --- Hoja1 (contagios acumulados por region) ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$row1 = @(43936, 44, 129, 66, 192, 13, 68, 330, 4682, 54, 167, 639, 559, 882, 153, 399, 7, 467, 8807)
for ($i = 0; $i -lt $row1.Length; $i++) {
    $ws1.Cells.Item(45, $i + 1).Value = $row1[$i]
}
$ws1.Range("A45").NumberFormat = "DD/MM/YY"

# --- Hoja2 (casos nuevos diarios por region) ---
$ws2 = $wb.Worksheets.Item("Hoja2")

$row2 = @(43936, 44, 1, 0, 1, 0, 0, 4, 48, 0, 4, 10, 2, 21, 3, 5, 0, 6, 105)
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws2.Cells.Item(45, $i + 1).Value = $row2[$i]
}
$ws2.Range("A45").NumberFormat = "DD/MM/YY"

# --- Hoja2 view: topLeftCell A31 -> A34, selection -> rows 46:60 plus N45, active cell N45 ---
$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 34
[void]$excel.Union($ws2.Range("46:60"), $ws2.Range("N45")).Select()
[void]$ws2.Range("N45").Activate()

# --- Hoja1 view: topLeftCell A34 -> A31, selection -> rows 46:60, active cell A60 ---
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 31
[void]$ws1.Range("46:60").Select()
[void]$ws1.Range("A60").Activate()
